$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 960, pushing existing rows 960:1035 down to 962:1037.
$ws.Rows("960:961").Insert()

# New row 960 (Calidad = Primera, Fecha serial 45013 = 2023-03-28)
$ws.Cells.Item(960, 1).Value2 = 3
$ws.Cells.Item(960, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(960, 3).Value2 = "Coquimbo"
$ws.Cells.Item(960, 4).Value2 = 45013
$ws.Cells.Item(960, 5).Value2 = 5
$ws.Cells.Item(960, 6).Value2 = 100114014
$ws.Cells.Item(960, 7).Value2 = "Betarraga"
$ws.Cells.Item(960, 8).Value2 = "Sin especificar"
$ws.Cells.Item(960, 9).Value2 = "Primera"
$ws.Cells.Item(960, 10).Value2 = 1600
$ws.Cells.Item(960, 11).Value2 = 750
$ws.Cells.Item(960, 12).Value2 = 750
$ws.Cells.Item(960, 13).Value2 = 750
$ws.Cells.Item(960, 14).Value2 = "`$/paquete 4 unidades"
$ws.Cells.Item(960, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(960, 16).Value2 = 188
$ws.Cells.Item(960, 17).Value2 = 4
$ws.Cells.Item(960, 18).Value2 = "Hortaliza"

# New row 961 (Calidad = Segunda, Fecha serial 45013 = 2023-03-28)
$ws.Cells.Item(961, 1).Value2 = 3
$ws.Cells.Item(961, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(961, 3).Value2 = "Coquimbo"
$ws.Cells.Item(961, 4).Value2 = 45013
$ws.Cells.Item(961, 5).Value2 = 5
$ws.Cells.Item(961, 6).Value2 = 100114014
$ws.Cells.Item(961, 7).Value2 = "Betarraga"
$ws.Cells.Item(961, 8).Value2 = "Sin especificar"
$ws.Cells.Item(961, 9).Value2 = "Segunda"
$ws.Cells.Item(961, 10).Value2 = 850
$ws.Cells.Item(961, 11).Value2 = 550
$ws.Cells.Item(961, 12).Value2 = 550
$ws.Cells.Item(961, 13).Value2 = 550
$ws.Cells.Item(961, 14).Value2 = "`$/paquete 4 unidades"
$ws.Cells.Item(961, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(961, 16).Value2 = 138
$ws.Cells.Item(961, 17).Value2 = 4
$ws.Cells.Item(961, 18).Value2 = "Hortaliza"

# Ensure the Fecha column keeps the same date/time number format as the rest of column D
$ws.Range("D960:D961").NumberFormat = $ws.Range("D962").NumberFormat

Write-Output "Done"
